# Split the sentence "... teste, imprime uma linha em branco." so that an
# en dash ("–") is inserted between the "i" and "mprime" of "imprime",
# turning the single run into three runs:
#   1) "    - Se não for o último caso de teste, i"
#   2) "–"
#   3) "mprime uma linha em branco."

$d = $word.ActiveDocument

# Locate the target sentence (use a duplicate range for the search so we
# don't disturb $d.Content itself).
$search = $d.Content.Duplicate
$found = $search.Find.Execute(
    "Se não for o último caso de teste, imprime uma linha em branco.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Target sentence not found"
}

# Find where "imprime" starts inside the matched text, then locate the
# point right after its leading "i".
$matchText = $search.Text
$offset = $matchText.IndexOf("imprime")
$splitPos = $search.Start + $offset + 1

# 1) Insert the en dash character right after the "i".
$insertionPoint = $d.Range($splitPos, $splitPos)
$insertionPoint.InsertAfter([string][char]0x2013)

# 2) Force the newly inserted dash into its own run (distinct from the
#    text before and after it) by nudging a character-formatting
#    property on just that character.
$dashRange = $d.Range($splitPos, $splitPos + 1)
$dashRange.Bold = 1
$dashRange.Bold = 0

Write-Host "Result: $($d.Range($search.Start, $search.End + 1).Text)"
